$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.50%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.03%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.161"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.17%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08383"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.30%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.948"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.23%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9731"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.03%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.72%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.53%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1893"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.31%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09706"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.23%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04614"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.99%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.16%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001294"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.57%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005727"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-6.03%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.401"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.83%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.444"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.42%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3364"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.86%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.647"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-15.22%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1363"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "'0.41%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04166"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.82%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001234"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-5.68%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004413"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.99%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001302"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.72%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-20.16%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02747"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.17%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05649"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.29%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007821"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.13%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1412"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.10%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007390"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.84%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002112"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.83%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007910"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.94%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3500"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006912"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.62%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003492"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.90%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003537"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'40.41%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.17%"
$ws.Range("E51").Style = "Normal"
